# dt_full_qoq_EXPORT_AVERAGE_1_9.xlsx
# Adds a new forecast-vintage column (BB) for "as-of" date 45986 (2025-11-25)
# alongside the existing BA column (45891 / 2025-08-22), and appends a new
# quarter row (83: 46934 / 2028-06-30) now that the forecast horizon extends
# one quarter further.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column BB values for rows 1-82 (row 1 is the header/as-of date row;
# rows 2-82 are the QoQ forecast values for that vintage).
$bbValues = @(
    45986, 1.63165117816655, 1.358419700277324, -0.1398451755242718, -0.07003068004286206,
    -7.3868590581191, -10.33903131837766, -1.060235412222937, 3.423433284297019, 2.063081011733999,
    2.99202665168275, 6.534207423404695, 1.294218807309846, 2.007778863461724, 2.488825190489734,
    -0.02236378853282872, 1.864148130528193, 0.3074416423962276, 1.591452302439862, 1.24229439238384,
    1.021051258256691, -2.083516682089652, 0.8945075486844729, 1.373553924371535, 0.3276711086097635,
    1.809509533486136, 0.4339813219297071, 0.4041977135476031, 1.507492882068462, 1.350099422002103,
    1.114800352984872, 1.516137977472326, -0.02409266101658147, -0.7388860519741201, 1.567695497950282,
    1.091382109524247, -0.2679292251141305, 1.682996656587392, 1.591260107191601, 1.13706115148689,
    1.735453665039003, 2.625533283765208, -0.3243298885145123, 0.6372341859553217, -0.5789332341234967,
    0.1728728569232914, 1.7, -1.4, 1.3, -0.6,
    -3.292009884772611, -20.353166912592, 17.87631681612835, 4.395171409529297, 1.763269038133103,
    1.651360288740889, -0.6023499639201475, 3.773751128807561, 0.3486937127494798, 0.6021692060546258,
    1.505090528591751, -1.274647749701984, 0.4647476577308112, -0.6885297541760451, -0.9272934801906558,
    -0.8572418502401149, 0.05364604092510206, 1.830682919206694, -2.604576060482884, -3.124437332092583,
    2.411315004676197, 0.2886831937783967, -0.7196185376451893, -0.7196185376451893, -0.7196185376451893,
    -0.7196185376451893, -0.7196185376451893, -0.7196185376451893, -0.7196185376451893, -0.7196185376451893,
    -0.7196185376451893, -0.7196185376451893
)

# Column BB is the 54th column (A=1 ... BA=53, BB=54).
$bbCol = 54

# Seed BB1's formatting from BA1 (date number format, font, border,
# alignment) since it is a header/date cell like the rest of row 1.
$ws.Range("BA1").Copy($ws.Range("BB1"))

for ($i = 0; $i -lt $bbValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, $bbCol).Value = $bbValues[$i]
}

# Append new row 83 for the next quarter (2028-06-30). Column A carries the
# same date style as the rest of column A; seed it from A82.
$ws.Range("A82").Copy($ws.Range("A83"))
$ws.Range("A83").Value = 46934
$ws.Range("BB83").Value = -0.7196185376451893
